$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update RF (column I) for rows 25 through 62 from 18.86026666666667 to 37.498
$ws.Range("I25:I62").Value = 37.498
